# Refined metadata to be additional tab
#
# 1) Refresh the "panel query" timestamps on the existing "data" sheet.
# 2) Add a new "metadata" sheet (after "data") summarising the panel
#    pull itself (name/id/version/request time/url).

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- 1. Update the per-row "time_taken" timestamps on the data sheet ----
$dataSheet.Range("F2").Value = "2021-10-05 14:23:04.916426"
$dataSheet.Range("F3").Value = "2021-10-05 14:23:04.916434"
$dataSheet.Range("F4").Value = "2021-10-05 14:23:04.916437"
$dataSheet.Range("F5").Value = "2021-10-05 14:23:04.916440"

# --- 2. Create the new "metadata" sheet, positioned after "data" -------
$newSheet = $wb.Worksheets.Add($null, $dataSheet)
$newSheet.Name = "metadata"

# Header row
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Reuse the bold/centered/bordered header style already used on "data"
# (style index 1) instead of fabricating a new one.
$dataSheet.Range("B1:F1").Copy()
$newSheet.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("F1").Copy()
$newSheet.Range("G1").PasteSpecial(-4122)

# Data row
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Upper gastrointestinal cancer pertinent cancer susceptibility"
$newSheet.Range("C2").Value = 273

# data_version must stay textual ("1.0"), not become the number 1.
# A text-formatted cell keeps a numeric-looking string as text instead
# of silently coercing it to a number; ClearFormats() afterwards drops
# the number-format override again (value/type already committed) so
# the cell ends up unstyled, matching the "data" sheet's own body cells.
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "1.0"
$newSheet.Range("D2").ClearFormats()

$newSheet.Range("E2").Value = "2017-11-05T02:37:20.339162Z"
$newSheet.Range("F2").Value = "2021-10-05 14:23:04.912855"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/273/?format=json"

# Same trick for A2's style (matches the index column styling on "data").
$dataSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

$newSheet.Range("A1").Select()

# Leave "data" as the active sheet/tab, as before.
$dataSheet.Activate()
